# Updated MCH102 to MCH251
# Add the new descriptive row (row 2) of collection metadata under the
# existing header row, matching the columns:
#   A: identifier, E: levelOfDescription, F: extentAndMedium, G: notes

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MCH113-1"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 21B | GRAP COUNT NUMER: NONE"

# Match the row's formatting (10pt Calibri, theme text color) across the
# whole row span that carries the new style, including the blank cells
# that flank the populated ones (C2, D2, H2).
$rowTwoCells = @("A2", "C2", "D2", "E2", "F2", "G2", "H2")
foreach ($addr in $rowTwoCells) {
    $cell = $ws.Range($addr)
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 10
    $cell.Font.ThemeColor = 1
}

# Keep the header row frozen (pane split after row 1) and move the active
# selection down to A3, just below the newly added data row.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A3").Select()
